$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.961.63"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "2.115.67"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").Value = "347.09"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").Value = "0.5185"
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").Value = "0.4443"
$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").Value = "53.83"
$ws.Range("E9").Value = "  +2.50%  "

$ws.Range("D10").Value = "0.09330"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").Value = "1.182"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").Value = "25.14"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").Value = "8.495"
$ws.Range("E13").Value = "  +3.75%  "

$ws.Range("D14").Value = "2.121.95"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D15").Value = "6.899"
$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("D16").Value = "103.15"
$ws.Range("E16").Value = "  +3.28%  "

$ws.Range("D17").Value = "0.00001162"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").Value = "21.48"
$ws.Range("E19").Value = "  +3.47%  "

$ws.Range("D20").Value = "0.06699"
$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").Value = "6.299"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").Value = "30.010.63"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "2.320"
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("D26").Value = "2.361.61"
$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("D27").Value = "22.10"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "2.543"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").Value = "162.71"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").Value = "134.15"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").Value = "1.146"
$ws.Range("E31").Value = "  -2.13%  "

$ws.Range("D32").Value = "1.776"
$ws.Range("E32").Value = "  +8.11%  "

$ws.Range("D33").Value = "0.1056"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "6.595"
$ws.Range("E35").Value = "  +6.11%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "3.977"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("E37").Value = "  +5.60%  "

$ws.Range("E38").Value = "  +1.76%  "

$ws.Range("E39").Value = "  +1.00%  "

$ws.Range("D40").Value = "0.7079"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("D41").Value = "12.70"
$ws.Range("E41").Value = "  +0.98%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.331"
$ws.Range("E42").Value = "  +1.10%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.2240"
$ws.Range("E43").Value = "  -2.38%  "

$ws.Range("D44").Value = "0.6847"

$ws.Range("D45").Value = "14.52"
$ws.Range("E45").Value = "  +1.56%  "

$ws.Range("D46").Value = "2.355"
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("D48").Value = "3.632"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "1.259"
$ws.Range("E49").Value = "  +8.48%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000356"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("E51").Value = "  -0.02%  "
